$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H26").Value = 0.65843
$ws.Range("I26").Value = 0.04403
$ws.Range("H27").Value = 0.12567
$ws.Range("I27").Value = 0.03588
$ws.Range("H28").Value = 0.66266
$ws.Range("I28").Value = 0.04577
$ws.Range("H29").Value = 0.18386
$ws.Range("I29").Value = 0.04795
$ws.Range("H30").Value = 0.65779
$ws.Range("I30").Value = 0.04485
$ws.Range("H31").Value = 0.12091
$ws.Range("I31").Value = 0.04084
$ws.Range("H32").Value = 0.66167
$ws.Range("I32").Value = 0.04663
$ws.Range("H33").Value = 0.1935
$ws.Range("I33").Value = 0.05546
$ws.Range("H34").Value = 0.663
$ws.Range("I34").Value = 0.03893
$ws.Range("H35").Value = 0.08211
$ws.Range("I35").Value = 0.04442
$ws.Range("H36").Value = 0.66238
$ws.Range("I36").Value = 0.03938
$ws.Range("H37").Value = 0.09186999999999999
$ws.Range("I37").Value = 0.04669
$ws.Range("H38").Value = 0.66398
$ws.Range("I38").Value = 0.04023
$ws.Range("H39").Value = 0.08211
$ws.Range("I39").Value = 0.04442
$ws.Range("H40").Value = 0.66269
$ws.Range("I40").Value = 0.03874
$ws.Range("H41").Value = 0.09186999999999999
$ws.Range("I41").Value = 0.04669
$ws.Range("H66").Value = 0.66328
$ws.Range("I66").Value = 0.04879
$ws.Range("H67").Value = 0.11603
$ws.Range("I67").Value = 0.056
$ws.Range("H68").Value = 0.66266
$ws.Range("I68").Value = 0.0467
$ws.Range("H69").Value = 0.12567
$ws.Range("I69").Value = 0.06551
$ws.Range("H70").Value = 0.66455
$ws.Range("I70").Value = 0.04736
$ws.Range("H71").Value = 0.11603
$ws.Range("I71").Value = 0.056
$ws.Range("H72").Value = 0.66651
$ws.Range("I72").Value = 0.04069
$ws.Range("H73").Value = 0.14019
$ws.Range("I73").Value = 0.05603
$ws.Range("H74").Value = 0.67132
$ws.Range("I74").Value = 0.03731
$ws.Range("H75").Value = 0.07247000000000001
$ws.Range("I75").Value = 0.05073
$ws.Range("H76").Value = 0.67194
$ws.Range("I76").Value = 0.03795
$ws.Range("H77").Value = 0.09651999999999999
$ws.Range("I77").Value = 0.05895
$ws.Range("H78").Value = 0.67196
$ws.Range("I78").Value = 0.03755
$ws.Range("H79").Value = 0.07247000000000001
$ws.Range("I79").Value = 0.05073
$ws.Range("H80").Value = 0.6722900000000001
$ws.Range("I80").Value = 0.0374
$ws.Range("H81").Value = 0.10139
$ws.Range("I81").Value = 0.05572
$ws.Range("H106").Value = 0.66351
$ws.Range("I106").Value = 0.03436
$ws.Range("H107").Value = 0.13055
$ws.Range("I107").Value = 0.04946
$ws.Range("H108").Value = 0.66325
$ws.Range("I108").Value = 0.04
$ws.Range("H109").Value = 0.16446
$ws.Range("I109").Value = 0.05392
$ws.Range("H110").Value = 0.66221
$ws.Range("I110").Value = 0.03539
$ws.Range("H111").Value = 0.1403
$ws.Range("I111").Value = 0.06009
$ws.Range("H112").Value = 0.66715
$ws.Range("I112").Value = 0.04361
$ws.Range("H113").Value = 0.17909
$ws.Range("I113").Value = 0.05892
$ws.Range("H114").Value = 0.67419
$ws.Range("I114").Value = 0.03071
$ws.Range("H115").Value = 0.09186999999999999
$ws.Range("I115").Value = 0.03567
$ws.Range("H116").Value = 0.67196
$ws.Range("I116").Value = 0.03191
$ws.Range("H117").Value = 0.1065
$ws.Range("I117").Value = 0.05034
$ws.Range("H118").Value = 0.67419
$ws.Range("I118").Value = 0.03071
$ws.Range("H119").Value = 0.09186999999999999
$ws.Range("I119").Value = 0.03567
$ws.Range("H120").Value = 0.67294
$ws.Range("I120").Value = 0.03176
$ws.Range("H121").Value = 0.11138
$ws.Range("I121").Value = 0.04545
$ws.Range("H146").Value = 0.66837
$ws.Range("I146").Value = 0.05208
$ws.Range("H147").Value = 0.15006
$ws.Range("I147").Value = 0.06028
$ws.Range("H148").Value = 0.6616300000000001
$ws.Range("I148").Value = 0.05656
$ws.Range("H149").Value = 0.16934
$ws.Range("I149").Value = 0.07485
$ws.Range("H150").Value = 0.66642
$ws.Range("I150").Value = 0.05443
$ws.Range("H151").Value = 0.15006
$ws.Range("I151").Value = 0.06028
$ws.Range("H152").Value = 0.6674600000000001
$ws.Range("I152").Value = 0.05491
$ws.Range("H153").Value = 0.20314
$ws.Range("I153").Value = 0.07555000000000001
$ws.Range("H154").Value = 0.67194
$ws.Range("I154").Value = 0.04406
$ws.Range("H155").Value = 0.11614
$ws.Range("I155").Value = 0.04911
$ws.Range("H156").Value = 0.67
$ws.Range("I156").Value = 0.04415
$ws.Range("H157").Value = 0.13554
$ws.Range("I157").Value = 0.06221
$ws.Range("H158").Value = 0.6712900000000001
$ws.Range("I158").Value = 0.04423
$ws.Range("H159").Value = 0.11614
$ws.Range("I159").Value = 0.04911
$ws.Range("H160").Value = 0.66839
$ws.Range("I160").Value = 0.04391
$ws.Range("H161").Value = 0.13554
$ws.Range("I161").Value = 0.06221
